$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Status moved from "Handed back: in sync with en-US" to the next stage of
# the handback pipeline. The 55b480a7 file is now back "In Translation"
# (Overview sheet only keeps the shared status text, cell refs untouched),
# while the fa4be397 file (and both per-locale table rows for both files)
# are now "Ready for handoff". Handback timestamps advance and an
# "out of date" Error Detail note is attached to each locale row.
# ---------------------------------------------------------------------------

$readyForHandoff = "Ready for handoff"
$msg55b480a7 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/923533a5f697198aab851c6136e03aea0e968096/e2e/55b480a7-9f5d-4e45-9679-e51c2e0c2521.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f80a0e9f7c5fec6759d712ddea0f457c44d6bc1/e2e/55b480a7-9f5d-4e45-9679-e51c2e0c2521.md."
$msgFa4be397 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/923533a5f697198aab851c6136e03aea0e968096/e2e/fa4be397-15a2-40e3-80d3-37cad9c27fac.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f80a0e9f7c5fec6759d712ddea0f457c44d6bc1/e2e/fa4be397-15a2-40e3-80d3-37cad9c27fac.md."

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 2 (55b480a7 file) stays on the same status text, which now reads
# "In Translation" instead of "Handed back: in sync with en-US".
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Row 3 (fa4be397 file) moves forward to "Ready for handoff".
$wsOverview.Range("E3").Value = $readyForHandoff
$wsOverview.Range("F3").Value = $readyForHandoff

# "Latest HO Xliff Generate Date" (shared by both rows) advances to the
# de-de handback timestamp.
$wsOverview.Range("G2").Value = "2016-11-08 23:40:38"
$wsOverview.Range("G3").Value = "2016-11-08 23:40:38"

# Columns E/F narrow from the old status-column width.
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $readyForHandoff
$wsZh.Range("H2").Value = "2016-11-08 23:40:25"
$wsZh.Range("P2").Value = $msg55b480a7

$wsZh.Range("C3").Value = $readyForHandoff
$wsZh.Range("H3").Value = "2016-11-08 23:40:25"
$wsZh.Range("P3").Value = $msgFa4be397

$wsZh.Columns.Item(3).ColumnWidth = 16.3333333333333
$wsZh.Columns.Item(16).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $readyForHandoff
$wsDe.Range("H2").Value = "2016-11-08 23:40:38"
$wsDe.Range("P2").Value = $msg55b480a7

$wsDe.Range("C3").Value = $readyForHandoff
$wsDe.Range("H3").Value = "2016-11-08 23:40:38"
$wsDe.Range("P3").Value = $msgFa4be397

$wsDe.Columns.Item(3).ColumnWidth = 16.3333333333333
$wsDe.Columns.Item(16).ColumnWidth = 39.1666666666667
